$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A56").NumberFormat = "@"
$ws.Range("A56").Value = "01-07-2021"
$ws.Range("A56").Style = "Normal"

$ws.Range("B56").Value = 3465060
$ws.Range("C56").Value = 31998
$ws.Range("D56").Value = 2668340
$ws.Range("E56").Value = 627955
$ws.Range("F56").Value = 168765
